$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: update title and link
$ws.Range("D6").Value = "Class 이해하기 :: Class를 쓰는 이유, Class vs function"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Class-%EC%9D%B4%ED%95%B4%ED%95%98%EA%B8%B0-Class%EB%A5%BC-%EC%93%B0%EB%8A%94-%EC%9D%B4%EC%9C%A0-Class-vs-function"

# Row 26: update title only
$ws.Range("D26").Value = "생성 모델의 새로운 흐름 확산 모델(Diffusion model)에 관하여"

# Row 28: update title and link
$ws.Range("D28").Value = "Robotarium - 원격 로봇 테스트 베드"
$ws.Range("E28").Value = "https://ropiens.tistory.com/185"

# Row 46: update title and link
$ws.Range("D46").Value = "백혈병 (Leukemia)"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/455"
